$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AF4").Value = 0.759
$ws.Range("AF5").Value = 0.944
$ws.Range("AF6").Value = 0.841
$ws.Range("AF7").Value = 0.9
$ws.Range("AF8").Value = 0.876
$ws.Range("AF9").Value = 0.778
$ws.Range("AF10").Value = 0.944
$ws.Range("AF11").Value = 0.944
$ws.Range("AF12").Value = 1.235
$ws.Range("AF13").Value = 1.444
